# Added ifoCAST full series evaluation:
# every existing data row (2..20) gains one more realised quarter of
# naive QoQ forecast error. Each row's values shift one column to the
# left (the oldest forecast-horizon value drops off the front) and a
# freshly computed error value is appended on the right - except for the
# rows that were already at the "triangle edge" (11..20), which simply
# lose their last value with nothing new appended (the series for that
# vintage has run out of realised data to compare against).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last populated column (as a 1-based column index, A=1) for each data
# row before this edit, and - where the row keeps its full width - the
# brand-new value to append in that last column.
$lastCol = @{
    2 = 11;  3 = 11;  4 = 11;  5 = 11;  6 = 11;
    7 = 11;  8 = 11;  9 = 11; 10 = 11; 11 = 11;
    12 = 10; 13 = 9; 14 = 8; 15 = 7; 16 = 6;
    17 = 5; 18 = 4; 19 = 3; 20 = 2
}

$appended = @{
    2  = -0.3644392301887736
    3  = 0.02750693478591659
    4  = -0.4322994165924858
    5  = 0.1131997290193177
    6  = 0.2163646915946629
    7  = 0.322776941072984
    8  = -0.4825338632108016
    9  = 0.4661714972207444
    10 = 1.11229800409388
}

for ($r = 2; $r -le 20; $r++) {
    $last = $lastCol[$r]

    # Read the current (pre-edit) row values first - column B (2) .. $last.
    $vals = @{}
    for ($c = 2; $c -le $last; $c++) {
        $vals[$c] = $ws.Cells.Item($r, $c).Value2
    }

    # Shift left: column c takes the old value that was in column c+1.
    for ($c = 2; $c -lt $last; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c + 1]
    }

    if ($appended.ContainsKey($r)) {
        # Row keeps its original width - append the newly realised value.
        $ws.Cells.Item($r, $last).Value = $appended[$r]
    }
    else {
        # Row shrinks by one column - drop the trailing cell entirely.
        $ws.Cells.Item($r, $last).ClearContents()
    }
}
